$d = $word.ActiveDocument

# Locate the "Remark." paragraph inside the Tests section. It currently reads
# " Remark.  " -- a plain leading space run, the bold/name-styled "Remark." +
# trailing space runs, and a final plain trailing space run. The edit removes
# the extraneous plain leading and trailing space runs, leaving just the
# name-styled "Remark." text (matching the pattern already used by the other
# similar paragraphs such as "Notation:", "Remarks:" and "Sketch Proof:").

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text
    $styleName = $p.Range.ParagraphStyle.NameLocal

    if ($styleName -eq "ProofStyle" -and $txt.Trim() -eq "Remark.") {
        $r = $p.Range
        $start = $r.Start
        $end = $r.End

        # Remove the trailing plain space run (the character immediately
        # before the paragraph mark).
        $trailing = $d.Range($end - 2, $end - 1)
        if ($trailing.Text -eq " ") {
            $trailing.Delete()
        }

        # Remove the leading plain space run.
        $leading = $d.Range($start, $start + 1)
        if ($leading.Text -eq " ") {
            $leading.Delete()
        }

        break
    }
}
